$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the NextAndBackFunctionality row: automated test cases 0 -> 4, status Testing -> Finished
$ws.Range("B4").Value = 4
$ws.Range("D4").Value = "Finished"

# Update the active selection to D5 as recorded in the saved view state
$ws.Range("D5").Select()
